$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: write a value into a cell as plain text, matching the original
# inline-string cell type, without leaving a residual NumberFormat style
# on the cell (some of these values - e.g. "314.94", "1.001" - look like
# numbers and would otherwise get auto-converted by Excel).
function Set-TextValue($cell, $value) {
    $cell.NumberFormat = "@"
    $cell.Value = $value
    $cell.Style = "Normal"
}

$updates = @{
    2  = @{ D = "27.611.63";     E = "  -2.09%  " }
    3  = @{ D = "1.842.95";      E = "  -1.18%  " }
    4  = @{ D = "1.001";         E = "  -0.42%  " }
    5  = @{ D = "314.94";        E = "  -1.16%  " }
    6  = @{ E = "  -0.27%  " }
    7  = @{ D = "0.4244" }
    8  = @{ D = "0.3643";        E = "  -1.10%  " }
    9  = @{ D = "45.29";         E = "  +0.49%  " }
    10 = @{ D = "0.07271";       E = "  -2.81%  " }
    11 = @{ D = "0.8928";        E = "  -4.33%  " }
    12 = @{ D = "20.62";         E = "  -3.00%  " }
    13 = @{ D = "1.827.34";      E = "  -7.18%  " }
    14 = @{ D = "6.570";         E = "  -2.02%  " }
    15 = @{ D = "5.359";         E = "  -1.37%  " }
    16 = @{ D = "0.06869";       E = "  -0.40%  " }
    17 = @{ E = "  -0.20%  " }
    18 = @{ E = "  -3.87%  " }
    19 = @{ D = "0.000008834";   E = "  -1.94%  " }
    20 = @{ D = "0.9997";        E = "  -0.38%  " }
    21 = @{ D = "15.47";         E = "  -2.40%  " }
    22 = @{ D = "27.587.04";     E = "  -2.13%  " }
    23 = @{ D = "4.987";         E = "  -2.35%  " }
    24 = @{ E = "  -1.94%  " }
    25 = @{ D = "2.047.61";      E = "  -6.06%  " }
    26 = @{ D = "2.023";         E = "  +0.35%  " }
    27 = @{ D = "154.95";        E = "  -0.11%  " }
    28 = @{ D = "18.52";         E = "  +0.90%  " }
    29 = @{ D = "5.220";         E = "  -1.72%  " }
    30 = @{ D = "118.80";        E = "  +5.15%  " }
    31 = @{ D = "1.834";         E = "  +6.65%  " }
    32 = @{ D = "0.08914";       E = "  -1.05%  " }
    33 = @{ D = "0.7782";        E = "  -1.52%  " }
    34 = @{ D = "4.570";         E = "  -5.37%  " }
    35 = @{ D = "2.956";         E = "  +0.44%  " }
    36 = @{ D = "1.102";         E = "  -5.85%  " }
    37 = @{ D = "0.9994";        E = "  -0.42%  " }
    38 = @{ D = "0.05408";       E = "  -0.30%  " }
    39 = @{ D = "1.098";         E = "  -2.24%  " }
    40 = @{ D = "0.01919";       E = "  -2.39%  " }
    41 = @{ D = "2.771";         E = "  -6.09%  " }
    42 = @{ D = "6.841";         E = "  -2.39%  " }
    43 = @{ D = "0.5062";        E = "  -3.49%  " }
    44 = @{ D = "0.1653";        E = "  -1.26%  " }
    45 = @{ D = "8.216";         E = "  -5.30%  " }
    46 = @{ D = "0.06608";       E = "  -1.93%  " }
    47 = @{ D = "10.34";         E = "  -2.04%  " }
    48 = @{ D = "0.4693";        E = "  -3.45%  " }
    49 = @{ D = "105.00";        E = "  -1.75%  " }
    50 = @{ D = "0.9995";        E = "  -0.30%  " }
    51 = @{ E = "  -2.59%  " }
}

foreach ($rowKey in $updates.Keys) {
    $rowData = $updates[$rowKey]
    if ($rowData.ContainsKey("D")) {
        Set-TextValue $ws.Range("D$rowKey") $rowData["D"]
    }
    if ($rowData.ContainsKey("E")) {
        Set-TextValue $ws.Range("E$rowKey") $rowData["E"]
    }
}
